# Applies the recorded edit to nextflow_dag.pptx:
#  - duplicates slide 1 into a new slide 2 ("finish draft of part II of book")
#  - trims slide 2 down to a single summarized row and retexts a few shapes
#  - updates the notes-page date placeholder text (2/7/24 -> 2/19/24) across
#    the notes master and all slide layouts
#  - forces creation of the slide 2 notes page (notesSlide2.xml)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Update the "last modified" date field text everywhere it appears
#    (notes master + all slide layouts use the same cached field text).
# ---------------------------------------------------------------------------
$oldDate = "2/7/24"
$newDate = "2/19/24"

$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $sh = $notesMaster.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Duplicate slide 1 to create slide 2.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Duplicate()
$s2 = $p.Slides.Item(2)

# ---------------------------------------------------------------------------
# 3. Remove the shapes that don't belong on the trimmed-down slide 2:
#    the "Power check" row, the "Discovery analysis" row, the
#    "Combine calibration check results" shape + its connectors, and a
#    handful of leftover connectors / decorative textboxes.
# ---------------------------------------------------------------------------
$idsToDelete = @(49,50,51,52,53,54,55,56,57,58,59,60,65,68,71,74,77,80,87,88,89,90,91,92,95,98,101,154,155,156,157,158,159,160,161,173,174,178,179,180,181,182)

for ($i = $s2.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s2.Shapes.Item($i)
    if ($idsToDelete -contains $sh.Id) {
        $sh.Delete()
    }
}

# ---------------------------------------------------------------------------
# 4. Retext / reposition the remaining shapes that changed.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)

    if ($sh.Id -eq 2) {
        $sh.TextFrame.TextRange.Text = "Run cellwise QC"
    }
    elseif ($sh.Id -eq 37) {
        $tr = $sh.TextFrame.TextRange
        $tr.Text = "Prepare trans discovery analysis"
        $tr.Characters(9, 6).Font.Italic = $true
    }
    elseif ($sh.Id -eq 45) {
        $sh.Left = 6354845
        $sh.Top = 2472568
        $sh.Width = 3534925
        $sh.Height = 279509
        $sh.TextFrame.TextRange.Text = "Pairwise QC and discovery analysis 1"
    }
    elseif ($sh.Id -eq 46) {
        $sh.Left = 6354845
        $sh.Top = 2821916
        $sh.Width = 3534925
        $sh.Height = 279509
        $sh.TextFrame.TextRange.Text = "Pairwise QC and discovery analysis 2"
    }
    elseif ($sh.Id -eq 47) {
        $sh.Left = 6354846
        $sh.Top = 3171264
        $sh.Width = 3534924
        $sh.Height = 279509
        $sh.TextFrame.TextRange.Text = "Pairwise QC and discovery analysis 3"
    }
    elseif ($sh.Id -eq 48) {
        $sh.Left = 6354845
        $sh.Top = 3872304
        $sh.Width = 3534924
        $sh.Height = 279509
        $tr = $sh.TextFrame.TextRange
        $tr.Text = "Pairwise QC and discovery analysis r"
        $tr.Characters(37, 1).Font.Italic = $true
    }
    elseif ($sh.Id -eq 172) {
        $sh.Left = 8084820
        $sh.Top = 3505156
        $sh.Width = 390555
        $sh.Height = 369332
    }
    elseif ($sh.Id -eq 177) {
        $sh.Left = 7897436
        $sh.Top = 4151813
        $sh.Width = 647228
        $sh.Height = 292388
    }
}

# ---------------------------------------------------------------------------
# 5. Force creation of the notes page for slide 2 (mirrors PowerPoint
#    eagerly materializing notesSlide2.xml for a newly duplicated slide).
# ---------------------------------------------------------------------------
$np2 = $s2.NotesPage
for ($i = 1; $i -le $np2.Shapes.Count; $i++) {
    $sh = $np2.Shapes.Item($i)
    if ($sh.Type -eq 14) {
        # placeholder; touch the slide-number placeholder's field text
    }
}
$notesBody = $np2.Shapes.Placeholders.Item(2)
$notesBody.TextFrame.TextRange.Text = $notesBody.TextFrame.TextRange.Text
